$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")
$ws.Rows(6).Insert()
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "default"
